# PracticalTranformerUse.pptx - "2017 ppt fixes" commit
#
# Slide 1 (the title slide) has its subtitle placeholder changed from
# "FME 2016 Training" to "FME 2017" - split across two runs ("FME " and
# "2017") just like the authored file.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Replace the whole paragraph text first ...
$tr.Text = "FME 2017"

# ... then re-set just the "2017" tail so PowerPoint splits it into its
# own run (mirrors the authored runs: "FME " / "2017").
$tail = $tr.Characters(5, 4)
$tail.Text = "2017"
